$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "Exp 11"
$ws.Range("B7").Value = 0.85
$ws.Range("C7").Value = 1
$ws.Range("F7").Value = "Exp11.png"

$ws.Range("F8").Select()
